$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(43, 8).Value = 3072.8262
$ws.Cells.Item(43, 9).Value = 810
$ws.Cells.Item(43, 10).Value = 4279.6665
$ws.Cells.Item(43, 11).Value = 810
$ws.Cells.Item(43, 12).Value = 4279.6665
$ws.Cells.Item(43, 13).Value = -741
$ws.Cells.Item(43, 14).Value = -4417.6665
$ws.Cells.Item(76, 8).Value = 3089274.2
$ws.Cells.Item(76, 9).Value = 3009.875
$ws.Cells.Item(76, 10).Value = 5558285.5
$ws.Cells.Item(76, 11).Value = 3009.875
$ws.Cells.Item(76, 12).Value = 5558285.5
$ws.Cells.Item(76, 13).Value = -2694.875
$ws.Cells.Item(76, 14).Value = -5558915.5
$ws.Cells.Item(79, 8).Value = 3089274.2
$ws.Cells.Item(79, 9).Value = 3009.875
$ws.Cells.Item(79, 10).Value = 5558285.5
$ws.Cells.Item(79, 11).Value = 3009.875
$ws.Cells.Item(79, 12).Value = 5558285.5
$ws.Cells.Item(79, 13).Value = -1917.875
$ws.Cells.Item(79, 14).Value = -5560469.5
$ws.Cells.Item(113, 8).Value = 58827244
$ws.Cells.Item(113, 10).Value = 5138.5
$ws.Cells.Item(113, 12).Value = 5138.5
$ws.Cells.Item(113, 14).Value = -11646.5
$ws.Cells.Item(126, 8).Value = 30780
$ws.Cells.Item(126, 10).Value = 30780
$ws.Cells.Item(126, 12).Value = 30780
$ws.Cells.Item(126, 14).Value = -40660
$ws.Cells.Item(137, 8).Value = 101765.77
$ws.Cells.Item(137, 9).Value = 161526.77
$ws.Cells.Item(137, 10).Value = 2164.1333
$ws.Cells.Item(137, 11).Value = 484580.3099999999
$ws.Cells.Item(137, 12).Value = 6492.3999
$ws.Cells.Item(137, 13).Value = -482030.3099999999
$ws.Cells.Item(137, 14).Value = -11592.3999
$ws.Cells.Item(138, 8).Value = 3956.8484
$ws.Cells.Item(138, 9).Value = 2860.4167
$ws.Cells.Item(138, 10).Value = 4583.381
$ws.Cells.Item(138, 11).Value = 8581.250100000001
$ws.Cells.Item(138, 12).Value = 13750.143
$ws.Cells.Item(138, 13).Value = -3441.250100000001
$ws.Cells.Item(138, 14).Value = -24030.143
$ws.Cells.Item(141, 8).Value = 1792.359
$ws.Cells.Item(141, 9).Value = 1662.3513
$ws.Cells.Item(141, 11).Value = 4987.0539
$ws.Cells.Item(141, 13).Value = 192.9461000000001

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 1109.4166
$ws.Cells.Item(2, 9).Value = 1114.375
$ws.Cells.Item(2, 10).Value = 1099.5
$ws.Cells.Item(2, 11).Value = 1114.375
$ws.Cells.Item(2, 12).Value = 1099.5
$ws.Cells.Item(2, 13).Value = -1001.375
$ws.Cells.Item(2, 14).Value = -1325.5
$ws.Cells.Item(4, 8).Value = 201
$ws.Cells.Item(4, 9).Value = 201
$ws.Cells.Item(4, 11).Value = 201
$ws.Cells.Item(4, 13).Value = -85
$ws.Cells.Item(32, 8).Value = 10240.431
$ws.Cells.Item(32, 9).Value = 7663.385
$ws.Cells.Item(32, 10).Value = 22205.285
$ws.Cells.Item(32, 11).Value = 7663.385
$ws.Cells.Item(32, 12).Value = 22205.285
$ws.Cells.Item(32, 13).Value = -7376.385
$ws.Cells.Item(32, 14).Value = -22779.285
$ws.Cells.Item(37, 8).Value = 36642.668
$ws.Cells.Item(37, 10).Value = 36642.668
$ws.Cells.Item(37, 12).Value = 36642.668
$ws.Cells.Item(37, 14).Value = -37188.668
$ws.Cells.Item(44, 8).Value = 39637.5
$ws.Cells.Item(44, 10).Value = 39637.5
$ws.Cells.Item(44, 12).Value = 39637.5
$ws.Cells.Item(44, 14).Value = -40613.5
$ws.Cells.Item(74, 8).Value = 32259728
$ws.Cells.Item(74, 10).Value = 3955.5557
$ws.Cells.Item(74, 12).Value = 3955.5557
$ws.Cells.Item(74, 14).Value = -5703.5557
$ws.Cells.Item(77, 8).Value = 32259728
$ws.Cells.Item(77, 10).Value = 3955.5557
$ws.Cells.Item(77, 12).Value = 19777.7785
$ws.Cells.Item(77, 14).Value = -28513.7785
$ws.Cells.Item(80, 8).Value = 36600
$ws.Cells.Item(80, 10).Value = 45000
$ws.Cells.Item(80, 12).Value = 45000
$ws.Cells.Item(80, 14).Value = -46996
$ws.Cells.Item(83, 8).Value = 36600
$ws.Cells.Item(83, 10).Value = 45000
$ws.Cells.Item(83, 12).Value = 135000
$ws.Cells.Item(83, 14).Value = -144984
$ws.Cells.Item(116, 8).Value = 1109.4166
$ws.Cells.Item(116, 9).Value = 1114.375
$ws.Cells.Item(116, 10).Value = 1099.5
$ws.Cells.Item(116, 11).Value = 1114.375
$ws.Cells.Item(116, 12).Value = 1099.5
$ws.Cells.Item(116, 13).Value = 1179.625
$ws.Cells.Item(116, 14).Value = -5687.5
$ws.Cells.Item(132, 8).Value = 11124478
$ws.Cells.Item(132, 9).Value = 13159908
$ws.Cells.Item(132, 10).Value = 75005.71000000001
$ws.Cells.Item(132, 11).Value = 39479724
$ws.Cells.Item(132, 12).Value = 225017.13
$ws.Cells.Item(132, 13).Value = -39477194
$ws.Cells.Item(132, 14).Value = -230077.13

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 1109.4166
$ws.Cells.Item(3, 9).Value = 1114.375
$ws.Cells.Item(3, 10).Value = 1099.5
$ws.Cells.Item(3, 11).Value = 1114.375
$ws.Cells.Item(3, 12).Value = 1099.5
$ws.Cells.Item(3, 13).Value = -1000.375
$ws.Cells.Item(3, 14).Value = -1327.5
$ws.Cells.Item(94, 8).Value = 1011.36
$ws.Cells.Item(94, 9).Value = 531.73334
$ws.Cells.Item(94, 10).Value = 1730.8
$ws.Cells.Item(94, 11).Value = 531.73334
$ws.Cells.Item(94, 12).Value = 1730.8
$ws.Cells.Item(94, 13).Value = -80.73334
$ws.Cells.Item(94, 14).Value = -2632.8
$ws.Cells.Item(134, 8).Value = 4699.4517
$ws.Cells.Item(134, 10).Value = 5500
$ws.Cells.Item(134, 12).Value = 16500
$ws.Cells.Item(134, 14).Value = -21570

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 6379.846
$ws.Cells.Item(31, 9).Value = 1530
$ws.Cells.Item(31, 10).Value = 6784
$ws.Cells.Item(31, 11).Value = 1530
$ws.Cells.Item(31, 12).Value = 6784
$ws.Cells.Item(31, 13).Value = -1235
$ws.Cells.Item(31, 14).Value = -7374
$ws.Cells.Item(34, 8).Value = 6379.846
$ws.Cells.Item(34, 9).Value = 1530
$ws.Cells.Item(34, 10).Value = 6784
$ws.Cells.Item(34, 11).Value = 1530
$ws.Cells.Item(34, 12).Value = 6784
$ws.Cells.Item(34, 13).Value = -1328
$ws.Cells.Item(34, 14).Value = -7188

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(69, 8).Value = 1880
$ws.Cells.Item(69, 10).Value = 1750
$ws.Cells.Item(69, 12).Value = 5250
$ws.Cells.Item(69, 14).Value = -6872
$ws.Cells.Item(70, 8).Value = 3983.4
$ws.Cells.Item(70, 9).Value = 1583.3334
$ws.Cells.Item(70, 11).Value = 4750.0002
$ws.Cells.Item(70, 13).Value = -4435.0002
$ws.Cells.Item(72, 8).Value = 1880
$ws.Cells.Item(72, 10).Value = 1750
$ws.Cells.Item(72, 12).Value = 15750
$ws.Cells.Item(72, 14).Value = -23862
$ws.Cells.Item(73, 8).Value = 3983.4
$ws.Cells.Item(73, 9).Value = 1583.3334
$ws.Cells.Item(73, 11).Value = 4750.0002
$ws.Cells.Item(73, 13).Value = -3658.0002
$ws.Cells.Item(82, 8).Value = 8340.833000000001
$ws.Cells.Item(82, 10).Value = 10011.25
$ws.Cells.Item(82, 12).Value = 30033.75
$ws.Cells.Item(82, 14).Value = -30845.75
$ws.Cells.Item(85, 8).Value = 8340.833000000001
$ws.Cells.Item(85, 10).Value = 10011.25
$ws.Cells.Item(85, 12).Value = 30033.75
$ws.Cells.Item(85, 14).Value = -32841.75
$ws.Cells.Item(131, 8).Value = 708.02
$ws.Cells.Item(131, 10).Value = 751.73865
$ws.Cells.Item(131, 12).Value = 2255.21595
$ws.Cells.Item(131, 14).Value = -12335.21595
$ws.Cells.Item(134, 8).Value = 5255
$ws.Cells.Item(134, 9).Value = 2047.9166
$ws.Cells.Item(134, 10).Value = 18083.334
$ws.Cells.Item(134, 11).Value = 6143.7498
$ws.Cells.Item(134, 12).Value = 54250.00199999999
$ws.Cells.Item(134, 13).Value = -1073.7498
$ws.Cells.Item(134, 14).Value = -64390.00199999999
$ws.Cells.Item(140, 8).Value = 3561.2666
$ws.Cells.Item(140, 9).Value = 1845.5714
$ws.Cells.Item(140, 11).Value = 5536.7142
$ws.Cells.Item(140, 13).Value = -356.7142000000003
$ws.Cells.Item(141, 8).Value = 5757.5
$ws.Cells.Item(141, 9).Value = 5757.5
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 17272.5
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -12092.5
$ws.Cells.Item(141, 14).ClearContents()

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 93
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(43, 8).Value = 672000
$ws.Cells.Item(43, 9).Value = 8000
$ws.Cells.Item(43, 10).Value = 2000000
$ws.Cells.Item(43, 11).Value = 8000
$ws.Cells.Item(43, 12).Value = 2000000
$ws.Cells.Item(43, 13).Value = -7849
$ws.Cells.Item(43, 14).Value = -2000302
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 5548393
$ws.Cells.Item(132, 10).Value = 103999.4
$ws.Cells.Item(132, 12).Value = 311998.2
$ws.Cells.Item(132, 14).Value = -317058.2

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(62, 8).Value = 4745.364
$ws.Cells.Item(62, 10).Value = 4966.5
$ws.Cells.Item(62, 12).Value = 4966.5
$ws.Cells.Item(62, 14).Value = -6214.5
$ws.Cells.Item(65, 8).Value = 4745.364
$ws.Cells.Item(65, 10).Value = 4966.5
$ws.Cells.Item(65, 12).Value = 24832.5
$ws.Cells.Item(65, 14).Value = -31072.5
